# Auto-generated: update crypto price/volume figures (columns D and E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.521.27"
$ws.Range("D3").Value = "2.612.07"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.545"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0818"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.58"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "3.005.54"
$ws.Range("D15").Value = "2.601.43"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.852"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "43.599.91"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "0.0₃0971"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "256.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0815"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.113"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0307"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "2.026.58"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "2.850.93"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.195"
$ws.Range("D51").Style = "Normal"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("E3").Value = "  +3.48%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("E6").Value = "  +3.34%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +2.79%  "
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("E13").Value = "  +3.18%  "
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("E15").Value = "  +3.24%  "
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("E18").Value = "  +1.69%  "
$ws.Range("E19").Value = "  +3.34%  "
$ws.Range("E20").Value = "  -1.78%  "
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("E24").Value = "  +1.60%  "
$ws.Range("E26").Value = "  +3.60%  "
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("E30").Value = "  +0.90%  "
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("E33").Value = "  +7.06%  "
$ws.Range("E34").Value = "  +3.65%  "
$ws.Range("E35").Value = "  +3.40%  "
$ws.Range("E36").Value = "  +3.61%  "
$ws.Range("E37").Value = "  -0.33%  "
$ws.Range("E38").Value = "  +9.59%  "
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("E41").Value = "  -3.91%  "
$ws.Range("E42").Value = "  +7.68%  "
$ws.Range("E43").Value = "  +0.71%  "
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("E47").Value = "  +3.01%  "
$ws.Range("E48").Value = "  -1.89%  "
$ws.Range("E49").Value = "  +3.04%  "
$ws.Range("E50").Value = "  +2.85%  "
$ws.Range("E51").Value = "  +3.10%  "
